# Apply updated crypto price/volume data per the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.527.92"
$ws.Range("E2").Value = "  -4.41%  "
$ws.Range("D3").Value = "2.194.15"
$ws.Range("E3").Value = "  -7.08%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'486.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.08%  "
$ws.Range("D6").Value = "'125.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.29%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "'0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.25%  "
$ws.Range("D9").Value = "2.214.75"
$ws.Range("E9").Value = "  -6.32%  "
$ws.Range("D10").Value = "'0.0922"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.30%  "
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").Value = "'4.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.73%  "
$ws.Range("D13").Value = "'0.314"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("D14").Value = "2.584.46"
$ws.Range("E14").Value = "  -7.06%  "
$ws.Range("D15").Value = "'21.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").Value = "53.449.64"
$ws.Range("E16").Value = "  -4.48%  "
$ws.Range("D17").Value = "'0.0000128"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("D18").Value = "2.203.17"
$ws.Range("E18").Value = "  -6.74%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'3.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'9.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.64%  "
$ws.Range("D21").Value = "'293.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.46%  "
$ws.Range("D22").Value = "'6.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.98%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'62.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.04%  "
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").Value = "2.294.52"
$ws.Range("E28").Value = "  -7.27%  "
$ws.Range("D29").Value = "'7.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").Value = "'165.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.04%  "
$ws.Range("D31").Value = "'1.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.80%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0₃0663"
$ws.Range("E33").Value = "  -6.65%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'0.993"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "'5.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").Value = "'17.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("D39").Value = "'0.825"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.09%  "
$ws.Range("D40").Value = "'35.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "'3.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.10%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "'1.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.99%  "
$ws.Range("D44").Value = "'3.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("D45").Value = "'124.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("D46").Value = "'4.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("D47").Value = "'0.0880"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("D48").Value = "'0.532"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.14%  "
$ws.Range("D49").Value = "'230.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("D50").Value = "'0.0470"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("D51").Value = "'0.0200"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.45%  "
